# Update gh-pages output data (想去人数 / "want to go" counts) on the
# 展览 (Exhibition), 演出 (Performance) and 全部类型 (All Types) sheets.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibition) sheet ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1691
$ws1.Range("F3").Value  = 9184
$ws1.Range("F4").Value  = 117
$ws1.Range("F6").Value  = 716
$ws1.Range("F7").Value  = 1382
$ws1.Range("F8").Value  = 208
$ws1.Range("F9").Value  = 65
$ws1.Range("F11").Value = 5962
$ws1.Range("F13").Value = 392
$ws1.Range("F15").Value = 4630
$ws1.Range("F16").Value = 18
$ws1.Range("F18").Value = 1153
$ws1.Range("F19").Value = 35
$ws1.Range("F20").Value = 344
$ws1.Range("F21").Value = 35
$ws1.Range("F22").Value = 2
$ws1.Range("F23").Value = 263
$ws1.Range("F24").Value = 16
$ws1.Range("F25").Value = 3040
$ws1.Range("F26").Value = 131

# --- 演出 (Performance) sheet ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 47

# --- 全部类型 (All Types) sheet ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1691
$ws4.Range("F3").Value  = 9184
$ws4.Range("F4").Value  = 117
$ws4.Range("F5").Value  = 47
$ws4.Range("F7").Value  = 716
$ws4.Range("F8").Value  = 1382
$ws4.Range("F9").Value  = 208
$ws4.Range("F10").Value = 65
$ws4.Range("F12").Value = 5962
$ws4.Range("F14").Value = 392
$ws4.Range("F16").Value = 4630
$ws4.Range("F17").Value = 18
$ws4.Range("F19").Value = 1153
$ws4.Range("F20").Value = 35
$ws4.Range("F21").Value = 344
$ws4.Range("F22").Value = 35
$ws4.Range("F23").Value = 2
$ws4.Range("F24").Value = 263
$ws4.Range("F25").Value = 16
$ws4.Range("F26").Value = 3040
$ws4.Range("F28").Value = 131
